$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update title text for year 2023 -> 2022
$ws.Range("A1").Value = "Total Demand for All Users customers for 2022 (Mthembanji)"

# Update monthly demand values in column B (rows 3-14)
$ws.Range("B3").Value = 402.0241910578001
$ws.Range("B4").Value = 362.843
$ws.Range("B5").Value = 393.2419999999999
$ws.Range("B6").Value = 363.646
$ws.Range("B7").Value = 443.7479999999999
$ws.Range("B8").Value = 325.4370000008
$ws.Range("B9").Value = 399.6180000000001
$ws.Range("B10").Value = 515.6840000000001
$ws.Range("B11").Value = 481.9270000001
$ws.Range("B12").Value = 487.7779999999999
$ws.Range("B13").Value = 467.098725
$ws.Range("B14").Value = 472.5559416667
